$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet lists overdue payment periods for a worker. A new period (2509)
# needs to be added as a new row, the totals (Valor mora / Cant. Periodos)
# need to reflect it, and the old "last row" box-closing border formatting
# must move from the old last row (19) to the newly appended one (20).
# ---------------------------------------------------------------------------

# 1) Insert a new row right below the current last data row (19). This pushes
#    the signature block (old rows 24/25) down to rows 25/26, exactly like
#    Excel's native "Insert Row" command.
$ws.Rows("20").Insert()

# 2) The inserted row starts out with generic/default formatting. Row 19
#    still carries the special "closing" border (bottom box). Copy that
#    formatting down onto the new row 20 first …
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122) # xlPasteFormats

# 3) … then restore row 19 back to the regular "middle" row formatting
#    (same as rows 16-18) now that it is no longer the last row.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false

# 4) The "Periodo Mora" column is now centered for every data row.
$ws.Range("E16:E20").HorizontalAlignment = -4108 # xlCenter

# 5) Populate the new row with the additional overdue period (2509) for the
#    same worker, matching the pattern of the previous rows.
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "21811641"
$ws.Range("D20").Value = "MARIA DUBERLINA ROLDAN UPEGUI"
$ws.Range("E20").Value = "2509"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

# 6) Refresh the summary figures: total overdue value and period count.
$ws.Range("E11").Value = 261924
$ws.Range("F13").Value = 5
